# Updates cryptos list values per data refresh (Fri Nov 17 06:52:09 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.420.62"
$ws.Range("E2").Value = "  -2.95%  "
$ws.Range("D3").Value = "1.985.80"
$ws.Range("E3").Value = "  -3.50%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'245.63"
$ws.Range("E5").Value = "  -3.12%  "
$ws.Range("D6").Value = "'0.629"
$ws.Range("E6").Value = "  -3.94%  "
$ws.Range("D7").Value = "'59.09"
$ws.Range("E7").Value = "  -12.77%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").Value = "'57.24"
$ws.Range("E10").Value = "  -4.24%  "
$ws.Range("E11").Value = "  +7.20%  "
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "'23.50"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("E14").Value = "  -6.92%  "
$ws.Range("D15").Value = "'14.03"
$ws.Range("D16").Value = "2.275.98"
$ws.Range("E16").Value = "  -3.50%  "
$ws.Range("D17").Value = "'5.47"
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("D18").Value = "1.981.84"
$ws.Range("E18").Value = "  -3.61%  "
$ws.Range("D19").Value = "36.288.70"
$ws.Range("E19").Value = "  -2.98%  "
$ws.Range("D20").Value = "'70.49"
$ws.Range("E20").Value = "  -4.37%  "
$ws.Range("D21").Value = "0.0₃0878"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'5.32"
$ws.Range("E22").Value = "  -3.20%  "
$ws.Range("D23").Value = "'234.27"
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -4.69%  "
$ws.Range("E26").Value = "  -6.02%  "
$ws.Range("D27").Value = "'9.97"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "'162.37"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.134"
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'19.92"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("D32").Value = "'1.18"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D33").Value = "'4.91"
$ws.Range("E33").Value = "  -6.55%  "
$ws.Range("D34").Value = "'0.0659"
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("E35").Value = "  -5.98%  "
$ws.Range("D36").Value = "'6.23"
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'2.26"
$ws.Range("E38").Value = "  -7.74%  "
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("E40").Value = "  -6.31%  "
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("D42").Value = "'0.0969"
$ws.Range("E42").Value = "  -6.27%  "
$ws.Range("E43").Value = "  -4.86%  "
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("E45").Value = "  -5.04%  "
$ws.Range("D46").Value = "'16.27"
$ws.Range("E46").Value = "  -7.12%  "
$ws.Range("D47").Value = "'92.59"
$ws.Range("E47").Value = "  -4.83%  "
$ws.Range("D48").Value = "'7.49"
$ws.Range("E48").Value = "  -5.85%  "
$ws.Range("D49").Value = "1.366.45"
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("E50").Value = "  -4.10%  "
$ws.Range("D51").Value = "'45.48"
$ws.Range("E51").Value = "  -2.76%  "
